$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 590.9
$ws.Cells.Item(2, 9).Value = 514.875
$ws.Cells.Item(2, 11).Value = 514.875
$ws.Cells.Item(2, 13).Value = -401.875
$ws.Cells.Item(9, 8).Value = 637993.1
$ws.Cells.Item(9, 9).Value = 20489
$ws.Cells.Item(9, 11).Value = 20489
$ws.Cells.Item(9, 13).Value = -20320
$ws.Cells.Item(12, 8).Value = 231.25
$ws.Cells.Item(12, 9).Value = 231.25
$ws.Cells.Item(12, 11).Value = 231.25
$ws.Cells.Item(12, 13).Value = -61.25
$ws.Cells.Item(16, 8).Value = 2758.3333
$ws.Cells.Item(16, 10).Value = 2758.3333
$ws.Cells.Item(16, 12).Value = 2758.3333
$ws.Cells.Item(16, 14).Value = -3218.3333
$ws.Cells.Item(18, 8).Value = 341.3
$ws.Cells.Item(18, 9).Value = 341.3
$ws.Cells.Item(18, 11).Value = 341.3
$ws.Cells.Item(18, 13).Value = -57.30000000000001
$ws.Cells.Item(19, 8).Value = 739.3871
$ws.Cells.Item(19, 9).Value = 695.95
$ws.Cells.Item(19, 10).Value = 818.36365
$ws.Cells.Item(19, 11).Value = 695.95
$ws.Cells.Item(19, 12).Value = 818.36365
$ws.Cells.Item(19, 13).Value = -520.95
$ws.Cells.Item(19, 14).Value = -1168.36365
$ws.Cells.Item(40, 8).Value = 2845
$ws.Cells.Item(40, 9).Value = 2431.25
$ws.Cells.Item(40, 10).Value = 4500
$ws.Cells.Item(40, 11).Value = 2431.25
$ws.Cells.Item(40, 12).Value = 4500
$ws.Cells.Item(40, 13).Value = -2256.25
$ws.Cells.Item(40, 14).Value = -4850
$ws.Cells.Item(43, 8).Value = 18174.26
$ws.Cells.Item(43, 9).Value = 3687.1428
$ws.Cells.Item(43, 10).Value = 68879.164
$ws.Cells.Item(43, 11).Value = 3687.1428
$ws.Cells.Item(43, 12).Value = 68879.164
$ws.Cells.Item(43, 13).Value = -3618.1428
$ws.Cells.Item(43, 14).Value = -69017.164
$ws.Cells.Item(58, 8).Value = 3787.875
$ws.Cells.Item(58, 9).Value = 1077.5
$ws.Cells.Item(58, 10).Value = 6498.25
$ws.Cells.Item(58, 11).Value = 3232.5
$ws.Cells.Item(58, 12).Value = 19494.75
$ws.Cells.Item(58, 13).Value = -3082.5
$ws.Cells.Item(58, 14).Value = -19794.75
$ws.Cells.Item(138, 8).Value = 2709.47
$ws.Cells.Item(138, 10).Value = 2746.0825
$ws.Cells.Item(138, 12).Value = 8238.247499999999
$ws.Cells.Item(138, 14).Value = -18518.2475

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1641.9656
$ws.Cells.Item(2, 10).Value = 2999
$ws.Cells.Item(2, 12).Value = 2999
$ws.Cells.Item(2, 14).Value = -3225
$ws.Cells.Item(32, 8).Value = 22926.521
$ws.Cells.Item(32, 9).Value = 10042.549
$ws.Cells.Item(32, 11).Value = 10042.549
$ws.Cells.Item(32, 13).Value = -9755.549000000001
$ws.Cells.Item(74, 8).Value = 1594.75
$ws.Cells.Item(74, 9).Value = 1493.6666
$ws.Cells.Item(74, 11).Value = 1493.6666
$ws.Cells.Item(74, 13).Value = -619.6666
$ws.Cells.Item(77, 8).Value = 1594.75
$ws.Cells.Item(77, 9).Value = 1493.6666
$ws.Cells.Item(77, 11).Value = 7468.333000000001
$ws.Cells.Item(77, 13).Value = -3100.333000000001
$ws.Cells.Item(116, 8).Value = 1641.9656
$ws.Cells.Item(116, 10).Value = 2999
$ws.Cells.Item(116, 12).Value = 2999
$ws.Cells.Item(116, 14).Value = -7587
$ws.Cells.Item(132, 8).Value = 3711.913
$ws.Cells.Item(132, 9).Value = 3367.1875
$ws.Cells.Item(132, 11).Value = 10101.5625
$ws.Cells.Item(132, 13).Value = -7571.5625

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1641.9656
$ws.Cells.Item(3, 10).Value = 2999
$ws.Cells.Item(3, 12).Value = 2999
$ws.Cells.Item(3, 14).Value = -3227
$ws.Cells.Item(57, 8).Value = 0
$ws.Cells.Item(57, 10).Value = 0
$ws.Cells.Item(57, 12).Value = 0
$ws.Cells.Item(57, 14).ClearContents()
$ws.Cells.Item(99, 8).Value = 2472
$ws.Cells.Item(99, 9).Value = 2525.4
$ws.Cells.Item(99, 11).Value = 2525.4
$ws.Cells.Item(99, 13).Value = -1027.4
$ws.Cells.Item(105, 8).Value = 3255.4285
$ws.Cells.Item(105, 9).Value = 2358.4
$ws.Cells.Item(105, 11).Value = 2358.4
$ws.Cells.Item(105, 13).Value = -611.4000000000001
$ws.Cells.Item(107, 8).Value = 3199.5
$ws.Cells.Item(107, 10).Value = 3399
$ws.Cells.Item(107, 12).Value = 3399
$ws.Cells.Item(107, 14).Value = -7239
$ws.Cells.Item(134, 8).Value = 3832.0625
$ws.Cells.Item(134, 9).Value = 3177.8125
$ws.Cells.Item(134, 10).Value = 4486.3125
$ws.Cells.Item(134, 11).Value = 9533.4375
$ws.Cells.Item(134, 12).Value = 13458.9375
$ws.Cells.Item(134, 13).Value = -6998.4375
$ws.Cells.Item(134, 14).Value = -18528.9375
$ws.Cells.Item(136, 8).Value = 0
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 20100116
$ws.Cells.Item(4, 9).Value = 22222350
$ws.Cells.Item(4, 10).Value = 1000000
$ws.Cells.Item(4, 11).Value = 22222350
$ws.Cells.Item(4, 12).Value = 1000000
$ws.Cells.Item(4, 13).Value = -22222238
$ws.Cells.Item(4, 14).Value = -1000224
$ws.Cells.Item(31, 8).Value = 19573.441
$ws.Cells.Item(31, 9).Value = 17498.5
$ws.Cells.Item(31, 10).Value = 19850.1
$ws.Cells.Item(31, 11).Value = 17498.5
$ws.Cells.Item(31, 12).Value = 19850.1
$ws.Cells.Item(31, 13).Value = -17203.5
$ws.Cells.Item(31, 14).Value = -20440.1
$ws.Cells.Item(34, 8).Value = 19573.441
$ws.Cells.Item(34, 9).Value = 17498.5
$ws.Cells.Item(34, 10).Value = 19850.1
$ws.Cells.Item(34, 11).Value = 17498.5
$ws.Cells.Item(34, 12).Value = 19850.1
$ws.Cells.Item(34, 13).Value = -17296.5
$ws.Cells.Item(34, 14).Value = -20254.1
$ws.Cells.Item(58, 8).Value = 4731.143
$ws.Cells.Item(58, 9).Value = 5291.143
$ws.Cells.Item(58, 11).Value = 5291.143
$ws.Cells.Item(58, 13).Value = -5088.143
$ws.Cells.Item(62, 8).Value = 9751.272000000001
$ws.Cells.Item(62, 9).Value = 9947.4
$ws.Cells.Item(62, 11).Value = 9947.4
$ws.Cells.Item(62, 13).Value = -9323.4
$ws.Cells.Item(65, 8).Value = 9751.272000000001
$ws.Cells.Item(65, 9).Value = 9947.4
$ws.Cells.Item(65, 11).Value = 49737
$ws.Cells.Item(65, 13).Value = -46617
$ws.Cells.Item(88, 8).Value = 29152
$ws.Cells.Item(88, 10).Value = 29152
$ws.Cells.Item(88, 12).Value = 29152
$ws.Cells.Item(88, 14).Value = -29964
$ws.Cells.Item(91, 8).Value = 29152
$ws.Cells.Item(91, 10).Value = 29152
$ws.Cells.Item(91, 12).Value = 29152
$ws.Cells.Item(91, 14).Value = -31960
$ws.Cells.Item(132, 8).Value = 4630
$ws.Cells.Item(132, 9).Value = 4050
$ws.Cells.Item(132, 10).Value = 6950
$ws.Cells.Item(132, 11).Value = 12150
$ws.Cells.Item(132, 12).Value = 20850
$ws.Cells.Item(132, 13).Value = -9620
$ws.Cells.Item(132, 14).Value = -25910
$ws.Cells.Item(136, 8).Value = 4731.143
$ws.Cells.Item(136, 9).Value = 5291.143
$ws.Cells.Item(136, 11).Value = 15873.429
$ws.Cells.Item(136, 13).Value = -13323.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 250
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 250
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 750
$ws.Cells.Item(5, 13).ClearContents()
$ws.Cells.Item(5, 14).Value = -974
$ws.Cells.Item(14, 8).Value = 306.4737
$ws.Cells.Item(14, 9).Value = 306.4737
$ws.Cells.Item(14, 11).Value = 919.4211
$ws.Cells.Item(14, 13).Value = -746.4211
$ws.Cells.Item(80, 8).Value = 2998.6667
$ws.Cells.Item(80, 10).Value = 2998.4
$ws.Cells.Item(80, 12).Value = 8995.200000000001
$ws.Cells.Item(80, 14).Value = -10867.2
$ws.Cells.Item(83, 8).Value = 2998.6667
$ws.Cells.Item(83, 10).Value = 2998.4
$ws.Cells.Item(83, 12).Value = 26985.6
$ws.Cells.Item(83, 14).Value = -36345.60000000001
$ws.Cells.Item(132, 8).Value = 2430.4194
$ws.Cells.Item(132, 9).Value = 1370.2858
$ws.Cells.Item(132, 11).Value = 12332.5722
$ws.Cells.Item(132, 13).Value = -9802.572200000001
$ws.Cells.Item(135, 8).Value = 250
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 10).Value = 250
$ws.Cells.Item(135, 11).Value = 0
$ws.Cells.Item(135, 12).Value = 2250
$ws.Cells.Item(135, 13).ClearContents()
$ws.Cells.Item(135, 14).Value = -7320

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(13, 8).Value = 164.66667
$ws.Cells.Item(13, 9).Value = 164.66667
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 11).Value = 164.66667
$ws.Cells.Item(13, 12).Value = 0
$ws.Cells.Item(13, 13).Value = -25.66667000000001
$ws.Cells.Item(13, 14).ClearContents()
$ws.Cells.Item(126, 8).Value = 3896.853
$ws.Cells.Item(126, 9).Value = 3348.5386
$ws.Cells.Item(126, 10).Value = 4236.2856
$ws.Cells.Item(126, 11).Value = 10045.6158
$ws.Cells.Item(126, 12).Value = 12708.8568
$ws.Cells.Item(126, 13).Value = -7575.6158
$ws.Cells.Item(126, 14).Value = -17648.8568
$ws.Cells.Item(132, 8).Value = 3680.5
$ws.Cells.Item(132, 9).Value = 4050.1904
$ws.Cells.Item(132, 11).Value = 12150.5712
$ws.Cells.Item(132, 13).Value = -9620.5712

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(21, 8).Value = 3949.5
$ws.Cells.Item(21, 9).Value = 4000
$ws.Cells.Item(21, 10).Value = 3899
$ws.Cells.Item(21, 11).Value = 4000
$ws.Cells.Item(21, 12).Value = 3899
$ws.Cells.Item(21, 13).Value = -3826
$ws.Cells.Item(21, 14).Value = -4247
$ws.Cells.Item(61, 8).Value = 6363
$ws.Cells.Item(61, 9).Value = 4908.4
$ws.Cells.Item(61, 11).Value = 4908.4
$ws.Cells.Item(61, 13).Value = -4706.4
$ws.Cells.Item(113, 8).Value = 6363
$ws.Cells.Item(113, 9).Value = 4908.4
$ws.Cells.Item(113, 11).Value = 4908.4
$ws.Cells.Item(113, 13).Value = -2738.4
$ws.Cells.Item(132, 8).Value = 3371.1538
$ws.Cells.Item(132, 9).Value = 2279.375
$ws.Cells.Item(132, 10).Value = 3856.389
$ws.Cells.Item(132, 11).Value = 6838.125
$ws.Cells.Item(132, 12).Value = 11569.167
$ws.Cells.Item(132, 13).Value = -4308.125
$ws.Cells.Item(132, 14).Value = -16629.167

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 1267.7273
$ws.Cells.Item(107, 9).Value = 993.8889
$ws.Cells.Item(107, 10).Value = 2500
$ws.Cells.Item(107, 11).Value = 2981.6667
$ws.Cells.Item(107, 12).Value = 7500
$ws.Cells.Item(107, 13).Value = -1061.6667
$ws.Cells.Item(107, 14).Value = -11340
